$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows before the current row 634, shifting existing
# rows 634:694 down to 636:696 (matches the dimension growing to A1:R696).
$ws.Rows("634:635").Insert()

# New row 634
$ws.Cells.Item(634, 1).Value = 5
$ws.Cells.Item(634, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(634, 3).Value = "Maule"
$ws.Cells.Item(634, 4).Value = 44769
$ws.Cells.Item(634, 5).Value = 7
$ws.Cells.Item(634, 6).Value = 100112020
$ws.Cells.Item(634, 7).Value = "Tomate"
$ws.Cells.Item(634, 8).Value = "Larga vida"
$ws.Cells.Item(634, 9).Value = "Primera"
$ws.Cells.Item(634, 10).Value = 2500
$ws.Cells.Item(634, 11).Value = 8500
$ws.Cells.Item(634, 12).Value = 8500
$ws.Cells.Item(634, 13).Value = 8500
$ws.Cells.Item(634, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(634, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(634, 16).Value = 472
$ws.Cells.Item(634, 17).Value = 18
$ws.Cells.Item(634, 18).Value = "Hortaliza"

# New row 635
$ws.Cells.Item(635, 1).Value = 5
$ws.Cells.Item(635, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(635, 3).Value = "Maule"
$ws.Cells.Item(635, 4).Value = 44769
$ws.Cells.Item(635, 5).Value = 7
$ws.Cells.Item(635, 6).Value = 100112020
$ws.Cells.Item(635, 7).Value = "Tomate"
$ws.Cells.Item(635, 8).Value = "Larga vida"
$ws.Cells.Item(635, 9).Value = "Primera"
$ws.Cells.Item(635, 10).Value = 1500
$ws.Cells.Item(635, 11).Value = 5500
$ws.Cells.Item(635, 12).Value = 5500
$ws.Cells.Item(635, 13).Value = 5500
$ws.Cells.Item(635, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(635, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(635, 16).Value = 550
$ws.Cells.Item(635, 17).Value = 10
$ws.Cells.Item(635, 18).Value = "Hortaliza"
